# Weekly update: insert the newest "Macroferia Regional de Talca - Mango" price
# record as a new row 78, shifting the existing historical rows (old 78-95) down
# to rows 79-96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 78 (pushes rows 78..95 down to 79..96,
# and extends the used range from A1:T95 to A1:T96).
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with this week's record.
$ws.Range("A78").Value = 5
$ws.Range("B78").Value = "Macroferia Regional de Talca"
$ws.Range("C78").Value = "Maule"
$ws.Range("D78").Value = 44543
$ws.Range("E78").Value = 7
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100108
$ws.Range("H78").Value = "Tropicales y subtropicales"
$ws.Range("I78").Value = 100108002
$ws.Range("J78").Value = "Mango"
$ws.Range("K78").Value = "Sin especificar"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 200
$ws.Range("N78").Value = 6000
$ws.Range("O78").Value = 6000
$ws.Range("P78").Value = 6000
$ws.Range("Q78").Value = "`$/bandeja 4 kilos"
$ws.Range("R78").Value = "Perú"
$ws.Range("S78").Value = 1500
$ws.Range("T78").Value = 4
